$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Coin name (B) and Link (C) cells that changed ---
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"

# --- Update Price (D) cells: force text format to preserve exact string (avoid numeric auto-conversion) ---
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }
$ws.Range("D2").Value = "30.279.34"
$ws.Range("D3").Value = "1.868.59"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "236.36"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.4671"
$ws.Range("D8").Value = "0.2837"
$ws.Range("D9").Value = "0.06544"
$ws.Range("D10").Value = "21.60"
$ws.Range("D11").Value = "0.07923"
$ws.Range("D12").Value = "98.07"
$ws.Range("D13").Value = "1.863.76"
$ws.Range("D14").Value = "5.151"
$ws.Range("D15").Value = "0.6795"
$ws.Range("D16").Value = "279.94"
$ws.Range("D17").Value = "30.290.96"
$ws.Range("D18").Value = "12.95"
$ws.Range("D19").Value = "1.0000"
$ws.Range("D20").Value = "5.450"
$ws.Range("D21").Value = "0.000007320"
$ws.Range("D22").Value = "2.112.69"
$ws.Range("D23").Value = "1.000"
$ws.Range("D24").Value = "6.163"
$ws.Range("D25").Value = "165.38"
$ws.Range("D26").Value = "9.158"
$ws.Range("D27").Value = "19.20"
$ws.Range("D28").Value = "1.945"
$ws.Range("D29").Value = "1.384"
$ws.Range("D30").Value = "0.09735"
$ws.Range("D31").Value = "4.397"
$ws.Range("D32").Value = "1.482"
$ws.Range("D33").Value = "4.131"
$ws.Range("D34").Value = "0.04713"
$ws.Range("D35").Value = "1.137"
$ws.Range("D36").Value = "0.7121"
$ws.Range("D37").Value = "2.724"
$ws.Range("D38").Value = "0.01863"
$ws.Range("D39").Value = "6.332"
$ws.Range("D40").Value = "2.532"
$ws.Range("D41").Value = "74.77"
$ws.Range("D42").Value = "1.970"
$ws.Range("D43").Value = "0.8540"
$ws.Range("D44").Value = "0.4197"
$ws.Range("D45").Value = "103.93"
$ws.Range("D46").Value = "0.9999"
$ws.Range("D47").Value = "7.234"
$ws.Range("D48").Value = "947.20"
$ws.Range("D49").Value = "9.255"
$ws.Range("D50").Value = "34.24"
$ws.Range("D51").Value = "0.1127"
foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# --- Update Volume(1h) (E) cells ---
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +6.42%  "
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("E15").Value = "  +2.15%  "
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("E40").Value = "  +1.66%  "
$ws.Range("E41").Value = "  +5.07%  "
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  -5.78%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  -0.60%  "
